$d = $word.ActiveDocument

# Locate the paragraph that currently ends the "没有任何的意义" sentence so we
# can insert the new sentence right after it, as a new BodyText paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    $text = $text.TrimEnd([char]13, [char]7)
    if ($text -eq "这里面写的东西没有任何的意义") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    $newPara = $target.Next()
    $newPara.Range.Text = "才怪，当然是为了测试"
    $newPara.Style = "BodyText"
}
